$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"1.336585149295133E-05"
$ws.Range("E2").Value = [double]"1.336585149295133E-05"

$ws.Range("D3").Value = 0.9847734688499893
$ws.Range("E3").Value = 0.9847734688499893

$ws.Range("D4").Value = [double]"8.363541502081753E-05"
$ws.Range("E4").Value = [double]"8.363541502081753E-05"

$ws.Range("D5").Value = [double]"3.715829035001071E-07"
$ws.Range("E5").Value = [double]"3.715829035001071E-07"

$ws.Range("D6").Value = 0.9482485481274607
$ws.Range("E6").Value = 0.9482485481274607

$ws.Range("C7").Value = $false
$ws.Range("D7").Value = 0.254417004346336
$ws.Range("E7").Value = 0.7455829956536639

$ws.Range("D8").Value = 0.9998962763499752
$ws.Range("E8").Value = [double]"0.0001037236500247563"

$ws.Range("D9").Value = 0.9953019302240896
$ws.Range("E9").Value = 0.004698069775910363

$ws.Range("D10").Value = 0.9999999999999978
$ws.Range("E10").Value = [double]"2.220446049250313E-15"

$ws.Range("D11").Value = 0.9999999912209664
$ws.Range("E11").Value = [double]"8.779033611183706E-09"
$ws.Range("F11").Value = 0.8519709706306458
$ws.Range("G11").Value = 0.7

$ws.Range("D12").Value = [double]"1.237259609185932E-05"
$ws.Range("E12").Value = [double]"1.237259609185932E-05"

$ws.Range("D13").Value = 0.9983241831121784
$ws.Range("E13").Value = 0.9983241831121784

$ws.Range("D14").Value = 0.9977326672394022
$ws.Range("E14").Value = 0.9977326672394022

$ws.Range("D15").Value = [double]"1.291473394686917E-07"
$ws.Range("E15").Value = [double]"1.291473394686917E-07"

$ws.Range("D16").Value = 0.1895625038770426
$ws.Range("E16").Value = 0.1895625038770426

$ws.Range("D17").Value = 0.5490134636165138
$ws.Range("E17").Value = 0.4509865363834862

$ws.Range("D18").Value = 0.9999997881397356
$ws.Range("E18").Value = [double]"2.118602644429757E-07"

$ws.Range("D19").Value = 0.9990472359004475
$ws.Range("E19").Value = 0.0009527640995524722

$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0

$ws.Range("D21").Value = 0.9999999999987512
$ws.Range("E21").Value = [double]"1.248778858098376E-12"
$ws.Range("F21").Value = 1.329138517379761
